$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.204.94'
$ws.Range("E2").Value = '  +0.98%  '
$ws.Range("D3").Value = '2.313.42'
$ws.Range("E3").Value = '  +0.41%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '''544.37'
$ws.Range("E5").Value = '  -0.47%  '
$ws.Range("D6").Value = '''131.83'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '''0.585'
$ws.Range("E8").Value = '  +2.32%  '
$ws.Range("D9").Value = '2.312.67'
$ws.Range("E9").Value = '  +0.48%  '
$ws.Range("D10").Value = '''0.100'
$ws.Range("E10").Value = '  -1.13%  '
$ws.Range("D11").Value = '''5.49'
$ws.Range("E11").Value = '  -0.02%  '
$ws.Range("E12").Value = '  +0.51%  '
$ws.Range("E13").Value = '  -0.08%  '
$ws.Range("D14").Value = '''23.83'
$ws.Range("E14").Value = '  +0.06%  '
$ws.Range("D15").Value = '2.726.69'
$ws.Range("E15").Value = '  +0.48%  '
$ws.Range("D16").Value = '59.091.07'
$ws.Range("E16").Value = '  +0.86%  '
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").Value = '2.313.41'
$ws.Range("E18").Value = '  -0.46%  '
$ws.Range("D19").Value = '''10.59'
$ws.Range("E19").Value = '  -0.58%  '
$ws.Range("D20").Value = '''4.19'
$ws.Range("E20").Value = '  -2.58%  '
$ws.Range("D21").Value = '''314.36'
$ws.Range("E21").Value = '  +0.18%  '
$ws.Range("E22").Value = '  +2.03%  '
$ws.Range("E23").Value = '  +0.18%  '
$ws.Range("D24").Value = '''62.66'
$ws.Range("E24").Value = '  -0.74%  '
$ws.Range("E25").Value = '  +3.41%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").Value = '''7.95'
$ws.Range("E27").Value = '  -0.90%  '
$ws.Range("D28").Value = '''1.31'
$ws.Range("E28").Value = '  +0.40%  '
$ws.Range("E29").Value = '  -0.73%  '
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").Value = '''171.12'
$ws.Range("E30").Value = '  +0.63%  '
$ws.Range("B31").Value = 'SuiNetwork'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D31").Value = '''1.18'
$ws.Range("E31").Value = '  +7.15%  '
$ws.Range("D32").Value = '0.0₃0736'
$ws.Range("E32").Value = '  +1.82%  '
$ws.Range("D33").Value = '''5.86'
$ws.Range("E33").Value = '  +1.83%  '
$ws.Range("E34").Value = '  +1.21%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '''1.33'
$ws.Range("E35").Value = '  +6.78%  '
$ws.Range("B36").Value = 'USDe'
$ws.Range("C36").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D36").Value = '''0.999'
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("B37").Value = 'EthereumClassic'
$ws.Range("C37").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D37").Value = '''17.83'
$ws.Range("E37").Value = '  +0.49%  '
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("D39").Value = '''4.07'
$ws.Range("E39").Value = '  +3.17%  '
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").Value = '''308.93'
$ws.Range("E40").Value = '  +4.46%  '
$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").Value = '''37.82'
$ws.Range("E41").Value = '  -0.75%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '''1.52'
$ws.Range("E42").Value = '  +0.60%  '
$ws.Range("D43").Value = '''141.21'
$ws.Range("E43").Value = '  +0.45%  '
$ws.Range("D44").Value = '''3.45'
$ws.Range("E44").Value = '  +0.61%  '
$ws.Range("D45").Value = '''0.0954'
$ws.Range("E45").Value = '  +0.49%  '
$ws.Range("D46").Value = '''0.0494'
$ws.Range("E46").Value = '  -1.23%  '
$ws.Range("D47").Value = '''0.558'
$ws.Range("E47").Value = '  +0.55%  '
$ws.Range("D48").Value = '''18.37'
$ws.Range("E48").Value = '  -0.19%  '
$ws.Range("E49").Value = '  -1.44%  '
$ws.Range("D50").Value = '''11.01'
$ws.Range("E51").Value = '  -0.35%  '
